$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: long Objetivos text is replaced in place by the "Docentes responsaveis" value (content moved) ---
$ws.Range("B10").Value = "5009972 - Gilberto Carvalho Coelho"
$ws.Range("C10").Value = "5009972 - Gilberto Carvalho Coelho"

# --- Rebuild rows 13-24 into rows 13-23 (content reshuffled, one row fewer) ---
$ws.Range("A13:C24").EntireRow.Delete()

# Row 13
$ws.Range("A8").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B8").Copy($ws.Range("B13"))
$ws.Range("B13").Value = "Semestral"
$ws.Range("C8").Copy($ws.Range("C13"))
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range("A8").Copy($ws.Range("A14"))
$ws.Range("A14").Value = "Short syllabus:"
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A8").Copy($ws.Range("A15"))
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C8").Copy($ws.Range("C15"))
$ws.Range("C15").Value = "01/01/2022"
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range("A8").Copy($ws.Range("A16"))
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A8").Copy($ws.Range("A17"))
$ws.Range("A17").Value = "Avaliação:"

# Row 18
$ws.Range("A8").Copy($ws.Range("A18"))
$ws.Range("A18").Value = "Método:"
$ws.Range("B8").Copy($ws.Range("B18"))
$ws.Range("B18").Value = "5009972 - Gilberto Carvalho Coelho"
$ws.Range("C8").Copy($ws.Range("C18"))
$ws.Range("C18").Value = "5009972 - Gilberto Carvalho Coelho"
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range("A8").Copy($ws.Range("A19"))
$ws.Range("A19").Value = "Critério:"
$ws.Range("B8").Copy($ws.Range("B19"))
$ws.Range("B19").Value = "O curso será ministrado na forma de aulas expositivas e aulas práticas em laboratório envolvendo preparação de amostras e caracterização microestrutural. Os resultados das aulas práticas serão apresentados oralmente e sujeitos a avaliação (T)."
$ws.Range("C8").Copy($ws.Range("C19"))
$ws.Range("C19").Value = "O curso será ministrado na forma de aulas expositivas e aulas práticas em laboratório envolvendo preparação de amostras e caracterização microestrutural. Os resultados das aulas práticas serão apresentados oralmente e sujeitos a avaliação (T)."
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A8").Copy($ws.Range("A20"))
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B8").Copy($ws.Range("B20"))
$ws.Range("B20").Value = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF) juntamente com a avaliação do trabalho prático (T). O critério para a nota final é:NF=((P1*0,8)+(T*0,2)+P2*1)/2"
$ws.Range("C8").Copy($ws.Range("C20"))
$ws.Range("C20").Value = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF) juntamente com a avaliação do trabalho prático (T). O critério para a nota final é:NF=((P1*0,8)+(T*0,2)+P2*1)/2"
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A8").Copy($ws.Range("A21"))
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B8").Copy($ws.Range("B21"))
$ws.Range("B21").Value = "Para os alunos que obtiverem 3,0≤NF<5,0, será aplicada uma avaliação de recuperação (R) que levará ao cálculo da média final (MF) com o seguinte critério:MF=(NF+R)/2"
$ws.Range("C8").Copy($ws.Range("C21"))
$ws.Range("C21").Value = "Para os alunos que obtiverem 3,0≤NF<5,0, será aplicada uma avaliação de recuperação (R) que levará ao cálculo da média final (MF) com o seguinte critério:MF=(NF+R)/2"
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Range("A8").Copy($ws.Range("A22"))
$ws.Range("A22").Value = "Requisitos:"

# Row 23
$ws.Range("B8").Copy($ws.Range("B23"))
$ws.Range("B23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Range("C8").Copy($ws.Range("C23"))
$ws.Range("C23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30
